$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3829411.6
$ws.Range("C7").Value = -13.81174274739063
$ws.Range("D7").Value = 3390
$ws.Range("E7").Value = 3390
$ws.Range("F7").Value = 1129.61994100295
$ws.Range("G7").Value = 20.40931750689026
